$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($range, [string]$text)
    if ($text -match '^[0-9]+(\.[0-9]+)?$') {
        # Purely numeric-looking text: force text storage so Excel does not
        # reinterpret it as a numeric value (preserving exact formatting,
        # leading/trailing zeros, etc.)
        $range.Value = "'" + $text
    } else {
        $range.Formula = $text
    }
}

Set-CellText $ws.Range("D2") '58.884.35'
Set-CellText $ws.Range("E2") '  -0.21%  '
Set-CellText $ws.Range("D3") '2.659.12'
Set-CellText $ws.Range("E3") '  +3.63%  '
Set-CellText $ws.Range("E4") '  +0.08%  '
Set-CellText $ws.Range("D5") '514.21'
Set-CellText $ws.Range("E5") '  +1.32%  '
Set-CellText $ws.Range("D6") '144.04'
Set-CellText $ws.Range("E6") '  +0.74%  '
Set-CellText $ws.Range("E7") '  -0.39%  '
Set-CellText $ws.Range("E8") '  +2.11%  '
Set-CellText $ws.Range("D9") '2.691.04'
Set-CellText $ws.Range("E9") '  +4.66%  '
Set-CellText $ws.Range("D10") '6.22'
Set-CellText $ws.Range("E10") '  +0.00%  '
Set-CellText $ws.Range("E11") '  +4.95%  '
Set-CellText $ws.Range("E12") '  +1.49%  '
Set-CellText $ws.Range("E13") '  -0.88%  '
Set-CellText $ws.Range("D14") '3.128.70'
Set-CellText $ws.Range("E14") '  +3.95%  '
Set-CellText $ws.Range("D15") '58.966.22'
Set-CellText $ws.Range("E15") '  +0.00%  '
Set-CellText $ws.Range("D16") '20.98'
Set-CellText $ws.Range("E16") '  +1.83%  '
Set-CellText $ws.Range("E17") '  +2.03%  '
Set-CellText $ws.Range("D18") '2.677.35'
Set-CellText $ws.Range("E18") '  +4.32%  '
Set-CellText $ws.Range("D19") '346.26'
Set-CellText $ws.Range("E19") '  +4.25%  '
Set-CellText $ws.Range("D20") '4.54'
Set-CellText $ws.Range("E20") '  +0.33%  '
Set-CellText $ws.Range("D21") '10.41'
Set-CellText $ws.Range("E21") '  +3.39%  '
Set-CellText $ws.Range("D22") '6.15'
Set-CellText $ws.Range("E22") '  +3.37%  '
Set-CellText $ws.Range("D23") '0.999'
Set-CellText $ws.Range("E23") '  +0.00%  '
Set-CellText $ws.Range("D24") '60.94'
Set-CellText $ws.Range("E24") '  +2.41%  '
Set-CellText $ws.Range("E25") '  +3.33%  '
Set-CellText $ws.Range("D26") '2.781.76'
Set-CellText $ws.Range("E26") '  +4.15%  '
Set-CellText $ws.Range("E27") '  -0.55%  '
Set-CellText $ws.Range("E28") '  +1.78%  '
Set-CellText $ws.Range("D29") '0.0₃0810'
Set-CellText $ws.Range("E29") '  +4.17%  '
Set-CellText $ws.Range("E30") '  +5.02%  '
Set-CellText $ws.Range("D31") '0.996'
Set-CellText $ws.Range("E31") '  -0.33%  '
Set-CellText $ws.Range("D32") '6.44'
Set-CellText $ws.Range("E32") '  +10.75%  '
Set-CellText $ws.Range("D33") '18.97'
Set-CellText $ws.Range("E33") '  +2.02%  '
Set-CellText $ws.Range("E34") '  +2.18%  '
Set-CellText $ws.Range("D35") '150.13'
Set-CellText $ws.Range("E35") '  +0.15%  '
Set-CellText $ws.Range("E36") '  +13.76%  '
Set-CellText $ws.Range("D37") '4.03'
Set-CellText $ws.Range("E37") '  +3.11%  '
Set-CellText $ws.Range("E38") '  +3.17%  '
Set-CellText $ws.Range("D39") '36.77'
Set-CellText $ws.Range("E39") '  +2.35%  '
Set-CellText $ws.Range("D40") '0.843'
Set-CellText $ws.Range("E40") '  +1.98%  '
Set-CellText $ws.Range("E41") '  +5.37%  '
Set-CellText $ws.Range("E42") '  +1.61%  '
Set-CellText $ws.Range("D43") '0.620'
Set-CellText $ws.Range("E43") '  +1.94%  '
Set-CellText $ws.Range("D44") '279.94'
Set-CellText $ws.Range("E44") '  -2.56%  '
Set-CellText $ws.Range("D45") '0.993'
Set-CellText $ws.Range("E45") '  -0.52%  '
Set-CellText $ws.Range("D46") '0.0981'
Set-CellText $ws.Range("E46") '  +0.14%  '
Set-CellText $ws.Range("D47") '19.82'
Set-CellText $ws.Range("E47") '  +6.49%  '
Set-CellText $ws.Range("D48") '0.0533'
Set-CellText $ws.Range("E48") '  +0.53%  '
Set-CellText $ws.Range("D49") '0.0231'
Set-CellText $ws.Range("E49") '  +1.84%  '
Set-CellText $ws.Range("B50") 'RenderToken'
Set-CellText $ws.Range("C50") 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-CellText $ws.Range("D50") '4.72'
Set-CellText $ws.Range("E50") '  +4.43%  '
Set-CellText $ws.Range("B51") 'Maker'
Set-CellText $ws.Range("C51") 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-CellText $ws.Range("D51") '2.001.87'
Set-CellText $ws.Range("E51") '  +4.68%  '
